# Add the "Distance_from_PFU_cm_correct" column (H) to the first worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mutliquant_pilot")

# Header
$ws.Range("H1").Value = "Distance_from_PFU_cm_correct"

# Data values for the new column
$values = @(0, 0.05, 0.1, 0.1, 0.15, 0.2, 0.2, 0.3, 0.4)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}

# Update the selection to match the post-edit state
$ws.Range("L18").Select()
